$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")

# Update ProductIdentifier (A), ProductName (B), GLN (E) and the duplicate
# ProductIdentifier (G) columns for the four product rows. The GS1 Company
# Prefix column (F) keeps its existing text values.
$ws.Range("A2").Value = 515381
$ws.Range("B2").Value = "Product14"
$ws.Range("E2").Value = 515381
$ws.Range("G2").Value = 515381

$ws.Range("A3").Value = 515382
$ws.Range("B3").Value = "Product15"
$ws.Range("E3").Value = 515382
$ws.Range("G3").Value = 515382

$ws.Range("A4").Value = 515383
$ws.Range("B4").Value = "Product16"
$ws.Range("E4").Value = 515383
$ws.Range("G4").Value = 515383

$ws.Range("A5").Value = 515384
$ws.Range("B5").Value = "Product17"
$ws.Range("E5").Value = 515384
$ws.Range("G5").Value = 515384

# Drop the stray row 10 leftover (E10) - shrinks the used range back to A1:L5.
$ws.Rows.Item(10).Delete()

# Make "Product" the active sheet/tab and restore the expected selection.
$ws.Activate()
$ws.Range("G2:G5").Select()
